$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("M2").Value = 1.533541666666667
$ws.Range("N2").Value = 4.600625
$ws.Range("O2").Value = 0.01998214594581092
$ws.Range("P2").Value = 0.01998214594581093
$ws.Range("Q2").Value = 0.0862857442361111
$ws.Range("R2").Value = 0.776571698125
$ws.Range("S2").Value = 0.005452383613742062
$ws.Range("T2").Value = 0.005452383613742063

$ws.Range("M3").Value = 3.948587333333334
$ws.Range("N3").Value = 11.845762
$ws.Range("O3").Value = 0.05145034536032411
$ws.Range("P3").Value = 0.05145034536032412
$ws.Range("Q3").Value = 0.2221698987015556
$ws.Range("R3").Value = 1.999529088314
$ws.Range("S3").Value = 0.01403888354758069
$ws.Range("T3").Value = 0.01403888354758069

$ws.Range("M4").Value = 70.69501233333334
$ws.Range("N4").Value = 212.085037
$ws.Range("O4").Value = 0.921160529766436
$ws.Range("P4").Value = 0.9211605297664361
$ws.Range("Q4").Value = 3.977701998943223
$ws.Range("R4").Value = 35.799317990489
$ws.Range("S4").Value = 0.2513504100983408
$ws.Range("T4").Value = 0.2513504100983408

$ws.Range("M5").Value = 0.568453
$ws.Range("N5").Value = 1.705359
$ws.Range("O5").Value = 0.007406978927428811
$ws.Range("P5").Value = 0.007406978927428812
$ws.Range("Q5").Value = 0.03198438701366667
$ws.Range("R5").Value = 0.287859483123
$ws.Range("S5").Value = 0.002021088757972569
$ws.Range("T5").Value = 0.00202108875797257

$ws.Range("M6").Value = 1.533541666666667
$ws.Range("N6").Value = 4.600625
$ws.Range("O6").Value = 0.01998214594581092
$ws.Range("P6").Value = 0.01998214594581093
$ws.Range("Q6").Value = 0.2299382151388889
$ws.Range("R6").Value = 2.06944393625
$ws.Range("S6").Value = 0.01452976233206886
$ws.Range("T6").Value = 0.01452976233206886

$ws.Range("M7").Value = 3.948587333333334
$ws.Range("N7").Value = 11.845762
$ws.Range("O7").Value = 0.05145034536032411
$ws.Range("P7").Value = 0.05145034536032412
$ws.Range("Q7").Value = 0.5920485523684446
$ws.Range("R7").Value = 5.328436971316001
$ws.Range("S7").Value = 0.03741146181274342
$ws.Range("T7").Value = 0.03741146181274343

$ws.Range("M8").Value = 70.69501233333334
$ws.Range("N8").Value = 212.085037
$ws.Range("O8").Value = 0.921160529766436
$ws.Range("P8").Value = 0.9211605297664361
$ws.Range("Q8").Value = 10.59996301925178
$ws.Range("R8").Value = 95.39966717326601
$ws.Range("S8").Value = 0.6698101196680952
$ws.Range("T8").Value = 0.6698101196680953

$ws.Range("M9").Value = 0.568453
$ws.Range("N9").Value = 1.705359
$ws.Range("O9").Value = 0.007406978927428811
$ws.Range("P9").Value = 0.007406978927428812
$ws.Range("Q9").Value = 0.08523346385133333
$ws.Range("R9").Value = 0.7671011746620001
$ws.Range("S9").Value = 0.005385890169456242
$ws.Range("T9").Value = 0.005385890169456242

